$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename "Intervention coverages" -> "Interventions coverages"
# ---------------------------------------------------------------------------
$covSheet = $wb.Worksheets.Item("Intervention coverages")
$covSheet.Name = "Interventions coverages"

# ---------------------------------------------------------------------------
# 2. Fill in the new intervention rows (4-10) on the coverages sheet
# ---------------------------------------------------------------------------
$covSheet.Range("A4").Value = "Complementary feeding 1"
$covSheet.Range("B4").Value = 0.0

$covSheet.Range("A5").Value = "Complementary feeding 2"
$covSheet.Range("B5").Value = 0.0

$covSheet.Range("A6").Value = "Complementary feeding 3"
$covSheet.Range("B6").Value = 0.0

$covSheet.Range("A7").Value = "Breastfeeding promotion"
$covSheet.Range("B7").Value = 0.0

$covSheet.Range("A8").Value = "IPTp"
$covSheet.Range("B8").Value = 0.0

$covSheet.Range("A9").Value = "BES"
$covSheet.Range("B9").Value = 0.0

$covSheet.Range("A10").Value = "MMS"
$covSheet.Range("B10").Value = 0.0

# ---------------------------------------------------------------------------
# 3. Add the three new worksheets (placed after the coverages sheet)
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$affected = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$affected.Name = "Interventions affected fraction"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$mortEff = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$mortEff.Name = "Interventions mortality eff"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$incEff = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$incEff.Name = "Interventions incidence eff"

# ---------------------------------------------------------------------------
# Helper to populate the common header row shared by the three new sheets
# ---------------------------------------------------------------------------
function Set-InterventionHeader($ws) {
    $ws.Range("A1").Value = "Interventions"
    $ws.Range("B1").Value = "Cause"
    $ws.Range("C1").Value = "<1 month"
    $ws.Range("D1").Value = "1-5 months"
    $ws.Range("E1").Value = "6-11 months"
    $ws.Range("F1").Value = "12-23 months"
    $ws.Range("G1").Value = "24-59 months"
}

function Set-InterventionBody($ws, $f2, $g2, $f3, $g3, $e4, $f4, $g4) {
    $ws.Range("A2").Value = "Zinc supplementation"
    $ws.Range("B2").Value = "Diarrhea"
    $ws.Range("C2").Value = 0.0
    $ws.Range("D2").Value = 0.0
    $ws.Range("E2").Value = 0.0
    $ws.Range("F2").Value = $f2
    $ws.Range("G2").Value = $g2

    $ws.Range("B3").Value = "Pneumonia"
    $ws.Range("C3").Value = 0.0
    $ws.Range("D3").Value = 0.0
    $ws.Range("E3").Value = 0.0
    $ws.Range("F3").Value = $f3
    $ws.Range("G3").Value = $g3

    $ws.Range("A4").Value = "Vitamin A supplementation"
    $ws.Range("B4").Value = "Diarrhea"
    $ws.Range("C4").Value = 0.0
    $ws.Range("D4").Value = 0.0
    $ws.Range("E4").Value = $e4
    $ws.Range("F4").Value = $f4
    $ws.Range("G4").Value = $g4

    $ws.Range("B5").Value = "Pneumonia"
    $ws.Range("C5").Value = 0.0
    $ws.Range("D5").Value = 0.0
    $ws.Range("E5").Value = 0.0
    $ws.Range("F5").Value = 0.0
    $ws.Range("G5").Value = 0.0
}

# ---------------------------------------------------------------------------
# 4. "Interventions affected fraction"
# ---------------------------------------------------------------------------
Set-InterventionHeader $affected
Set-InterventionBody $affected 0.253 0.253 0.253 0.253 0.416 0.416 0.416

# ---------------------------------------------------------------------------
# 5. "Interventions mortality eff"
# ---------------------------------------------------------------------------
Set-InterventionHeader $mortEff
Set-InterventionBody $mortEff 0.5 0.5 0.51 0.51 0.3 0.3 0.3

# ---------------------------------------------------------------------------
# 6. "Interventions incidence eff"
# ---------------------------------------------------------------------------
Set-InterventionHeader $incEff
Set-InterventionBody $incEff 0.65 0.65 0.52 0.52 0.62 0.62 0.62
